$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E7"  = 16.899
    "D9"  = -7.545
    "E12" = 17.646
    "E14" = 17.007
    "D18" = -8.543000000000001
    "D20" = -7.44
    "E26" = 16.938
    "D27" = -7.866000000000001
    "E27" = 16.777
    "E29" = 16.941
    "D35" = -7.606
    "E37" = 16.792
    "E38" = 16.835
    "E51" = 16.65
    "E52" = 16.657
    "E55" = 16.697
    "D69" = -7.679
    "E69" = 17.255
    "E70" = 17.609
    "D76" = -7.794
    "D78" = -8.123999999999999
    "E81" = 16.457
    "D82" = -8.398999999999999
    "D83" = -8.031000000000001
    "E83" = 16.767
    "D93" = -7.512
    "E102" = 16.701
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
